$d = $word.ActiveDocument

# ------------------------------------------------------------------
# The commit adds a short intro about what is needed to run the
# program on Linux. As part of that edit, the title paragraph (the
# very first paragraph of the document) is given bold / dark-red /
# larger formatting, and the document's "_GoBack" bookmark (which used
# to sit after the very last picture in the document) ends up instead
# sitting in the middle of that title's text - right after "...correr e"
# and before "l programa...". That's simply where the author's cursor
# was when they last edited, which Word tracks via _GoBack.
# ------------------------------------------------------------------

# 1) Remove the old "_GoBack" bookmark (currently at the end of the
#    document, right after the final inline picture) BEFORE creating
#    the new one, so the name lookup below can't grab the wrong one.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete() | Out-Null
}

# 2) Rebuild the title paragraph: bold, bold-complex-script, dark red
#    (C00000), size 28 half-points (14pt) on both the run text and the
#    paragraph mark, with the text split into two runs so the
#    "_GoBack" bookmark can sit between them.
$titleRange = $d.Paragraphs(1).Range
$titleXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:pPr><w:rPr><w:b/><w:bCs/><w:color w:val="C00000"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:b/><w:bCs/><w:color w:val="C00000"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr>' +
    '<w:t>¿Qué necesitamos para correr e</w:t></w:r>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
    '<w:r><w:rPr><w:b/><w:bCs/><w:color w:val="C00000"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr>' +
    '<w:t>l programa en un sistema Linux?</w:t></w:r>' +
    '</w:p>'
$titleRange.InsertXML($titleXml) | Out-Null

Write-Output "Title paragraph reformatted and _GoBack bookmark relocated."
